$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(375).Insert()

$ws.Range("A375").Value = 4
$ws.Range("B375").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C375").Value = "Los Lagos"
$ws.Range("D375").Value = 45093
$ws.Range("D375").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E375").Value = 10
$ws.Range("F375").Value = 100112040
$ws.Range("G375").Value = "Cilantro"
$ws.Range("H375").Value = "Sin especificar"
$ws.Range("I375").Value = "Primera"
$ws.Range("J375").Value = 200
$ws.Range("K375").Value = 12000
$ws.Range("L375").Value = 14000
$ws.Range("M375").Value = 13000
$ws.Range("N375").Value = "$/caja 36 atados"
$ws.Range("O375").Value = "Región Metropolitana"
$ws.Range("P375").Value = 361
$ws.Range("Q375").Value = 36
$ws.Range("R375").Value = "Hortaliza"
